# Update Angpt1-Itga5.xlsx worksheet with refreshed TPM-derived NATMI output.
# The old data had 9 data rows (Sending cluster in {ECs, FAPs, MuSCs} x
# Target cluster in {ECs, FAPs, MuSCs}). The new data drops the "ECs" sending
# cluster entirely (leaving 6 data rows: FAPs/MuSCs x ECs/FAPs/MuSCs) and
# refreshes all the numeric expression/specificity columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose "Sending cluster" was MuSCs under the old
# shared-string ordering (old rows 8,9,10) -- that cluster's data is gone
# from the refreshed export entirely, and the remaining rows shift so the
# final sheet has data rows 2..7 (dimension A1:T7).
$ws.Rows("8:10").Delete()

# Row 2: FAPs -> Angpt1 -> Itga5 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Angpt1"
$ws.Range("C2").Value = "Itga5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.95459633333333
$ws.Range("H2").Value = 38.863789
$ws.Range("I2").Value = 0.8906505749177925
$ws.Range("J2").Value = 0.8906505749177924
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 144.4028903198314
$ws.Range("R2").Value = 1299.626012878482
$ws.Range("S2").Value = 0.2311065532004074
$ws.Range("T2").Value = 0.2311065532004074

# Row 3: FAPs -> Angpt1 -> Itga5 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Angpt1"
$ws.Range("C3").Value = "Itga5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.95459633333333
$ws.Range("H3").Value = 38.863789
$ws.Range("I3").Value = 0.8906505749177925
$ws.Range("J3").Value = 0.8906505749177924
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.05649099999999
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 358.6544379227109
$ws.Range("R3").Value = 3227.889941304399
$ws.Range("S3").Value = 0.5740009133803606
$ws.Range("T3").Value = 0.5740009133803605

# Row 4: FAPs -> Angpt1 -> Itga5 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Angpt1"
$ws.Range("C4").Value = "Itga5"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.95459633333333
$ws.Range("H4").Value = 38.863789
$ws.Range("I4").Value = 0.8906505749177925
$ws.Range("J4").Value = 0.8906505749177924
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.125957666666667
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 53.45011606008856
$ws.Range("R4").Value = 481.051044540797
$ws.Range("S4").Value = 0.08554310833702455
$ws.Range("T4").Value = 0.08554310833702454

# Row 5: MuSCs -> Angpt1 -> Itga5 -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Angpt1"
$ws.Range("C5").Value = "Itga5"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.590497666666667
$ws.Range("H5").Value = 4.771493
$ws.Range("I5").Value = 0.1093494250822076
$ws.Range("J5").Value = 0.1093494250822076
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 17.72903255369267
$ws.Range("R5").Value = 159.561292983234
$ws.Range("S5").Value = 0.02837405536680615
$ws.Range("T5").Value = 0.02837405536680615

# Row 6: MuSCs -> Angpt1 -> Itga5 -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Angpt1"
$ws.Range("C6").Value = "Itga5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.590497666666667
$ws.Range("H6").Value = 4.771493
$ws.Range("I6").Value = 0.1093494250822076
$ws.Range("J6").Value = 0.1093494250822076
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.05649099999999
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("Q6").Value = 44.033718379007
$ws.Range("R6").Value = 396.303465411063
$ws.Range("S6").Value = 0.07047283372673717
$ws.Range("T6").Value = 0.07047283372673717

# Row 7: MuSCs -> Angpt1 -> Itga5 -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Angpt1"
$ws.Range("C7").Value = "Itga5"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.590497666666667
$ws.Range("H7").Value = 4.771493
$ws.Range("I7").Value = 0.1093494250822076
$ws.Range("J7").Value = 0.1093494250822076
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.125957666666667
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("Q7").Value = 6.562326041598779
$ws.Range("R7").Value = 59.06093437438901
$ws.Range("S7").Value = 0.01050253598866426
$ws.Range("T7").Value = 0.01050253598866426
